$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (sheet 1)
#   - New client row for "FABIMP BENIGNO BRAVO S.A.S." is inserted at row 9
#     (pushing "FRANK FERRETERIA FRANKFERRE CIA." down to row 9... actually
#     the FRANK row now carries the old FABIMP-era M8 value of 5372.02, and
#     the original row 8 (FRANK) is renamed to FABIMP with M8 reset to 0).
#   - "VIEJO RIVAS MAYRA ANABELLE" shifts from row 9 to row 10.
#   - The trailing "x de 8" summary row shifts from row 10 to row 11 and the
#     counters change from "de 8" to "de 9" (one more client row now).
# -----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a fresh row before the old row 9, shifting VIEJO RIVAS and the
# totals row down by one. The new row inherits the number formatting of
# the row above it (style 2), which is what the inserted client row needs.
$ws1.Rows.Item(9).Insert()

# Row 8 becomes the new client "FABIMP BENIGNO BRAVO S.A.S." with all-zero
# figures for this sheet.
$ws1.Range("B8").Value = "FABIMP BENIGNO BRAVO S.A.S."
$ws1.Cells.Item(8, 13).Value = 0

# Row 9 (brand new) becomes "FRANK FERRETERIA FRANKFERRE CIA." carrying the
# PORCELANATO (column M) figure that used to sit on row 8.
$ws1.Range("A9").Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws1.Range("B9").Value = "FRANK FERRETERIA FRANKFERRE CIA."
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(9, $col).Value = 0
}
$ws1.Cells.Item(9, 13).Value = 5372.02

# Row 11 is the shifted-down summary row; refresh its "x de 8" -> "x de 9"
# labels (still text cells, style 3 already carried over by the insert).
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(11, $col).Value = "0 de 9"
}
$ws1.Cells.Item(11, 13).Value = "2 de 9"

# -----------------------------------------------------------------------
# Sheet "VENTA MENSUAL" (sheet 2) - same client reshuffle, different
# monthly columns (C..G).
# -----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(9).Insert()

# Row 8 becomes "FABIMP BENIGNO BRAVO S.A.S." with its own monthly figures.
$ws2.Range("B8").Value = "FABIMP BENIGNO BRAVO S.A.S."
$ws2.Cells.Item(8, 3).Value = 1187.62
$ws2.Cells.Item(8, 4).Value = 0
$ws2.Cells.Item(8, 5).Value = 0
$ws2.Cells.Item(8, 6).Value = 0
$ws2.Cells.Item(8, 7).Value = 1000

# Row 9 (brand new) becomes "FRANK FERRETERIA FRANKFERRE CIA." carrying the
# septiembre (column F) figure that used to sit on row 8.
$ws2.Range("A9").Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws2.Range("B9").Value = "FRANK FERRETERIA FRANKFERRE CIA."
$ws2.Cells.Item(9, 3).Value = 0
$ws2.Cells.Item(9, 4).Value = 0
$ws2.Cells.Item(9, 5).Value = 0
$ws2.Cells.Item(9, 6).Value = 5372.02
$ws2.Cells.Item(9, 7).Value = 0

# Row 11 is the shifted-down totals row; update it for the new FABIMP
# contributions (presupuesto +1000, junio +1187.62).
$ws2.Cells.Item(11, 3).Value = 1187.62
$ws2.Cells.Item(11, 4).Value = 0
$ws2.Cells.Item(11, 5).Value = 0
$ws2.Cells.Item(11, 6).Value = 5415.120000000001
$ws2.Cells.Item(11, 7).Value = 1200

# Column C widens slightly (11 -> 13 raw OOXML units) to fit the new value.
$ws2.Columns.Item(3).ColumnWidth = 12.2
